# The "type" column (E) for every feature row (rows 4-8) is being
# standardized from the old "cont"/"cat" labels to a single "numeric"
# label.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4:E8").Value = "numeric"

# The saved cursor/selection in the sheet moves to E12.
$ws.Range("E12").Select()
